$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel stores them as text (matching original inlineStr text cells)
# instead of silently converting to a floating point number.
$textCells = @("D5","D6","D8","D9","D14","D18","D20","D21","D26","D27","D28","D30","D31","D33","D35","D39","D40","D41","D45","D46","D47","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '66.530.24'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '3.595.52'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '609.17'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').Value = '148.61'
$ws.Range('E6').Value = '  +2.68%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '0.489'
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('D9').Value = '8.08'
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').Value = '4.211.85'
$ws.Range('E12').Value = '  +1.32%  '
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('D14').Value = '29.83'
$ws.Range('E14').Value = '  -0.87%  '
$ws.Range('D15').Value = '3.579.27'
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('D16').Value = '66.664.55'
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('E17').Value = '  +0.80%  '
$ws.Range('D18').Value = '11.48'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('E19').Value = '  +3.05%  '
$ws.Range('D20').Value = '15.10'
$ws.Range('E20').Value = '  +1.44%  '
$ws.Range('D21').Value = '428.01'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').Value = '3.742.16'
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = '0.0000122'
$ws.Range('E26').Value = '  +3.41%  '
$ws.Range('D27').Value = '8.32'
$ws.Range('E27').Value = '  +4.02%  '
$ws.Range('D28').Value = '9.53'
$ws.Range('E28').Value = '  +4.76%  '
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').Value = '1.47'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').Value = '3.595.24'
$ws.Range('E32').Value = '  +1.28%  '
$ws.Range('D33').Value = '0.157'
$ws.Range('E33').Value = '  +1.96%  '
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('D35').Value = '7.86'
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('E38').Value = '  -2.43%  '
$ws.Range('D39').Value = '177.26'
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').Value = '0.0857'
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('D41').Value = '5.25'
$ws.Range('E41').Value = '  +0.75%  '
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('E44').Value = '  +9.15%  '
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '24.96'
$ws.Range('E46').Value = '  -3.31%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').Value = '1.17'
$ws.Range('E47').Value = '  -2.48%  '
$ws.Range('D48').Value = '24.02'
$ws.Range('E48').Value = '  +2.36%  '
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('D50').Value = '0.953'
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('D51').Value = '2.427.42'
$ws.Range('E51').Value = '  +5.35%  '
